$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 52.47848103381103

$ws.Range("N2:N6").Value = $newValue
